$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A187:A199").NumberFormat = "@"

$ws.Cells.Item(187, 1).Value = '2026-01-28'
$ws.Cells.Item(187, 2).Value = '16:26:03'
$ws.Cells.Item(187, 3).Value = '16:00'
$ws.Cells.Item(187, 4).Value = 'Bathroom'
$ws.Cells.Item(187, 5).Value = 'No Motion'
$ws.Cells.Item(187, 6).Value = 'Inactive'

$ws.Cells.Item(188, 1).Value = '2026-01-28'
$ws.Cells.Item(188, 2).Value = '16:26:05'
$ws.Cells.Item(188, 3).Value = '16:00'
$ws.Cells.Item(188, 4).Value = 'Bathroom'
$ws.Cells.Item(188, 5).Value = 'No Motion'
$ws.Cells.Item(188, 6).Value = 'Inactive'

$ws.Cells.Item(189, 1).Value = '2026-01-28'
$ws.Cells.Item(189, 2).Value = '16:26:08'
$ws.Cells.Item(189, 3).Value = '16:00'
$ws.Cells.Item(189, 4).Value = 'Bathroom'
$ws.Cells.Item(189, 5).Value = 'No Motion'
$ws.Cells.Item(189, 6).Value = 'Inactive'

$ws.Cells.Item(190, 1).Value = '2026-01-28'
$ws.Cells.Item(190, 2).Value = '16:26:13'
$ws.Cells.Item(190, 3).Value = '16:00'
$ws.Cells.Item(190, 4).Value = 'Bathroom'
$ws.Cells.Item(190, 5).Value = 'No Motion'
$ws.Cells.Item(190, 6).Value = 'Inactive'

$ws.Cells.Item(191, 1).Value = '2026-01-28'
$ws.Cells.Item(191, 2).Value = '16:26:18'
$ws.Cells.Item(191, 3).Value = '16:00'
$ws.Cells.Item(191, 4).Value = 'Bathroom'
$ws.Cells.Item(191, 5).Value = 'No Motion'
$ws.Cells.Item(191, 6).Value = 'Inactive'

$ws.Cells.Item(192, 1).Value = '2026-01-28'
$ws.Cells.Item(192, 2).Value = '16:26:23'
$ws.Cells.Item(192, 3).Value = '16:00'
$ws.Cells.Item(192, 4).Value = 'Bathroom'
$ws.Cells.Item(192, 5).Value = 'No Motion'
$ws.Cells.Item(192, 6).Value = 'Inactive'

$ws.Cells.Item(193, 1).Value = '2026-01-28'
$ws.Cells.Item(193, 2).Value = '16:26:28'
$ws.Cells.Item(193, 3).Value = '16:00'
$ws.Cells.Item(193, 4).Value = 'Bathroom'
$ws.Cells.Item(193, 5).Value = 'No Motion'
$ws.Cells.Item(193, 6).Value = 'Inactive'

$ws.Cells.Item(194, 1).Value = '2026-01-28'
$ws.Cells.Item(194, 2).Value = '16:26:33'
$ws.Cells.Item(194, 3).Value = '16:00'
$ws.Cells.Item(194, 4).Value = 'Bathroom'
$ws.Cells.Item(194, 5).Value = 'No Motion'
$ws.Cells.Item(194, 6).Value = 'Inactive'

$ws.Cells.Item(195, 1).Value = '2026-01-28'
$ws.Cells.Item(195, 2).Value = '16:26:38'
$ws.Cells.Item(195, 3).Value = '16:00'
$ws.Cells.Item(195, 4).Value = 'Bathroom'
$ws.Cells.Item(195, 5).Value = 'No Motion'
$ws.Cells.Item(195, 6).Value = 'Inactive'

$ws.Cells.Item(196, 1).Value = '2026-01-28'
$ws.Cells.Item(196, 2).Value = '16:26:43'
$ws.Cells.Item(196, 3).Value = '16:00'
$ws.Cells.Item(196, 4).Value = 'Bathroom'
$ws.Cells.Item(196, 5).Value = 'No Motion'
$ws.Cells.Item(196, 6).Value = 'Inactive'

$ws.Cells.Item(197, 1).Value = '2026-01-28'
$ws.Cells.Item(197, 2).Value = '16:26:48'
$ws.Cells.Item(197, 3).Value = '16:00'
$ws.Cells.Item(197, 4).Value = 'Bathroom'
$ws.Cells.Item(197, 5).Value = 'No Motion'
$ws.Cells.Item(197, 6).Value = 'Inactive'

$ws.Cells.Item(198, 1).Value = '2026-01-28'
$ws.Cells.Item(198, 2).Value = '16:26:53'
$ws.Cells.Item(198, 3).Value = '16:00'
$ws.Cells.Item(198, 4).Value = 'Bathroom'
$ws.Cells.Item(198, 5).Value = 'No Motion'
$ws.Cells.Item(198, 6).Value = 'Inactive'

$ws.Cells.Item(199, 1).Value = '2026-01-28'
$ws.Cells.Item(199, 2).Value = '16:26:58'
$ws.Cells.Item(199, 3).Value = '16:00'
$ws.Cells.Item(199, 4).Value = 'Bathroom'
$ws.Cells.Item(199, 5).Value = 'No Motion'
$ws.Cells.Item(199, 6).Value = 'Inactive'

$ws.Range("A187:A199").Style = "Normal"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A185:A198").NumberFormat = "@"
$ws.Range("E185:E198").NumberFormat = "@"

$ws.Cells.Item(185, 1).Value = '2026-01-28'
$ws.Cells.Item(185, 2).Value = '16:26:03'
$ws.Cells.Item(185, 3).Value = '16:00'
$ws.Cells.Item(185, 4).Value = 'Bathroom'
$ws.Cells.Item(185, 5).Value = '88.3%'
$ws.Cells.Item(185, 6).Value = 'Active'

$ws.Cells.Item(186, 1).Value = '2026-01-28'
$ws.Cells.Item(186, 2).Value = '16:26:04'
$ws.Cells.Item(186, 3).Value = '16:00'
$ws.Cells.Item(186, 4).Value = 'Bathroom'
$ws.Cells.Item(186, 5).Value = '86.9%'
$ws.Cells.Item(186, 6).Value = 'Active'

$ws.Cells.Item(187, 1).Value = '2026-01-28'
$ws.Cells.Item(187, 2).Value = '16:26:06'
$ws.Cells.Item(187, 3).Value = '16:00'
$ws.Cells.Item(187, 4).Value = 'Bathroom'
$ws.Cells.Item(187, 5).Value = '87.4%'
$ws.Cells.Item(187, 6).Value = 'Active'

$ws.Cells.Item(188, 1).Value = '2026-01-28'
$ws.Cells.Item(188, 2).Value = '16:26:10'
$ws.Cells.Item(188, 3).Value = '16:00'
$ws.Cells.Item(188, 4).Value = 'Bathroom'
$ws.Cells.Item(188, 5).Value = '88.3%'
$ws.Cells.Item(188, 6).Value = 'Active'

$ws.Cells.Item(189, 1).Value = '2026-01-28'
$ws.Cells.Item(189, 2).Value = '16:26:14'
$ws.Cells.Item(189, 3).Value = '16:00'
$ws.Cells.Item(189, 4).Value = 'Bathroom'
$ws.Cells.Item(189, 5).Value = '88.3%'
$ws.Cells.Item(189, 6).Value = 'Active'

$ws.Cells.Item(190, 1).Value = '2026-01-28'
$ws.Cells.Item(190, 2).Value = '16:26:18'
$ws.Cells.Item(190, 3).Value = '16:00'
$ws.Cells.Item(190, 4).Value = 'Bathroom'
$ws.Cells.Item(190, 5).Value = '88.3%'
$ws.Cells.Item(190, 6).Value = 'Active'

$ws.Cells.Item(191, 1).Value = '2026-01-28'
$ws.Cells.Item(191, 2).Value = '16:26:22'
$ws.Cells.Item(191, 3).Value = '16:00'
$ws.Cells.Item(191, 4).Value = 'Bathroom'
$ws.Cells.Item(191, 5).Value = '88.3%'
$ws.Cells.Item(191, 6).Value = 'Active'

$ws.Cells.Item(192, 1).Value = '2026-01-28'
$ws.Cells.Item(192, 2).Value = '16:26:30'
$ws.Cells.Item(192, 3).Value = '16:00'
$ws.Cells.Item(192, 4).Value = 'Bathroom'
$ws.Cells.Item(192, 5).Value = '88.4%'
$ws.Cells.Item(192, 6).Value = 'Active'

$ws.Cells.Item(193, 1).Value = '2026-01-28'
$ws.Cells.Item(193, 2).Value = '16:26:34'
$ws.Cells.Item(193, 3).Value = '16:00'
$ws.Cells.Item(193, 4).Value = 'Bathroom'
$ws.Cells.Item(193, 5).Value = '88.3%'
$ws.Cells.Item(193, 6).Value = 'Active'

$ws.Cells.Item(194, 1).Value = '2026-01-28'
$ws.Cells.Item(194, 2).Value = '16:26:39'
$ws.Cells.Item(194, 3).Value = '16:00'
$ws.Cells.Item(194, 4).Value = 'Bathroom'
$ws.Cells.Item(194, 5).Value = '87.4%'
$ws.Cells.Item(194, 6).Value = 'Active'

$ws.Cells.Item(195, 1).Value = '2026-01-28'
$ws.Cells.Item(195, 2).Value = '16:26:42'
$ws.Cells.Item(195, 3).Value = '16:00'
$ws.Cells.Item(195, 4).Value = 'Bathroom'
$ws.Cells.Item(195, 5).Value = '88.3%'
$ws.Cells.Item(195, 6).Value = 'Active'

$ws.Cells.Item(196, 1).Value = '2026-01-28'
$ws.Cells.Item(196, 2).Value = '16:26:50'
$ws.Cells.Item(196, 3).Value = '16:00'
$ws.Cells.Item(196, 4).Value = 'Bathroom'
$ws.Cells.Item(196, 5).Value = '88.3%'
$ws.Cells.Item(196, 6).Value = 'Active'

$ws.Cells.Item(197, 1).Value = '2026-01-28'
$ws.Cells.Item(197, 2).Value = '16:26:55'
$ws.Cells.Item(197, 3).Value = '16:00'
$ws.Cells.Item(197, 4).Value = 'Bathroom'
$ws.Cells.Item(197, 5).Value = '88.3%'
$ws.Cells.Item(197, 6).Value = 'Active'

$ws.Cells.Item(198, 1).Value = '2026-01-28'
$ws.Cells.Item(198, 2).Value = '16:27:02'
$ws.Cells.Item(198, 3).Value = '16:00'
$ws.Cells.Item(198, 4).Value = 'Bathroom'
$ws.Cells.Item(198, 5).Value = '88.3%'
$ws.Cells.Item(198, 6).Value = 'Active'

$ws.Range("A185:A198").Style = "Normal"
$ws.Range("E185:E198").Style = "Normal"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A185:A197").NumberFormat = "@"

$ws.Cells.Item(185, 1).Value = '2026-01-28'
$ws.Cells.Item(185, 2).Value = '16:26:04'
$ws.Cells.Item(185, 3).Value = '16:00'
$ws.Cells.Item(185, 4).Value = 'Bathroom'
$ws.Cells.Item(185, 5).Value = '22.8C'
$ws.Cells.Item(185, 6).Value = 'Active'

$ws.Cells.Item(186, 1).Value = '2026-01-28'
$ws.Cells.Item(186, 2).Value = '16:26:05'
$ws.Cells.Item(186, 3).Value = '16:00'
$ws.Cells.Item(186, 4).Value = 'Bathroom'
$ws.Cells.Item(186, 5).Value = '22.8C'
$ws.Cells.Item(186, 6).Value = 'Active'

$ws.Cells.Item(187, 1).Value = '2026-01-28'
$ws.Cells.Item(187, 2).Value = '16:26:07'
$ws.Cells.Item(187, 3).Value = '16:00'
$ws.Cells.Item(187, 4).Value = 'Bathroom'
$ws.Cells.Item(187, 5).Value = '22.8C'
$ws.Cells.Item(187, 6).Value = 'Active'

$ws.Cells.Item(188, 1).Value = '2026-01-28'
$ws.Cells.Item(188, 2).Value = '16:26:11'
$ws.Cells.Item(188, 3).Value = '16:00'
$ws.Cells.Item(188, 4).Value = 'Bathroom'
$ws.Cells.Item(188, 5).Value = '22.8C'
$ws.Cells.Item(188, 6).Value = 'Active'

$ws.Cells.Item(189, 1).Value = '2026-01-28'
$ws.Cells.Item(189, 2).Value = '16:26:15'
$ws.Cells.Item(189, 3).Value = '16:00'
$ws.Cells.Item(189, 4).Value = 'Bathroom'
$ws.Cells.Item(189, 5).Value = '22.8C'
$ws.Cells.Item(189, 6).Value = 'Active'

$ws.Cells.Item(190, 1).Value = '2026-01-28'
$ws.Cells.Item(190, 2).Value = '16:26:19'
$ws.Cells.Item(190, 3).Value = '16:00'
$ws.Cells.Item(190, 4).Value = 'Bathroom'
$ws.Cells.Item(190, 5).Value = '22.8C'
$ws.Cells.Item(190, 6).Value = 'Active'

$ws.Cells.Item(191, 1).Value = '2026-01-28'
$ws.Cells.Item(191, 2).Value = '16:26:23'
$ws.Cells.Item(191, 3).Value = '16:00'
$ws.Cells.Item(191, 4).Value = 'Bathroom'
$ws.Cells.Item(191, 5).Value = '22.7C'
$ws.Cells.Item(191, 6).Value = 'Active'

$ws.Cells.Item(192, 1).Value = '2026-01-28'
$ws.Cells.Item(192, 2).Value = '16:26:31'
$ws.Cells.Item(192, 3).Value = '16:00'
$ws.Cells.Item(192, 4).Value = 'Bathroom'
$ws.Cells.Item(192, 5).Value = '22.8C'
$ws.Cells.Item(192, 6).Value = 'Active'

$ws.Cells.Item(193, 1).Value = '2026-01-28'
$ws.Cells.Item(193, 2).Value = '16:26:35'
$ws.Cells.Item(193, 3).Value = '16:00'
$ws.Cells.Item(193, 4).Value = 'Bathroom'
$ws.Cells.Item(193, 5).Value = '22.8C'
$ws.Cells.Item(193, 6).Value = 'Active'

$ws.Cells.Item(194, 1).Value = '2026-01-28'
$ws.Cells.Item(194, 2).Value = '16:26:39'
$ws.Cells.Item(194, 3).Value = '16:00'
$ws.Cells.Item(194, 4).Value = 'Bathroom'
$ws.Cells.Item(194, 5).Value = '22.8C'
$ws.Cells.Item(194, 6).Value = 'Active'

$ws.Cells.Item(195, 1).Value = '2026-01-28'
$ws.Cells.Item(195, 2).Value = '16:26:43'
$ws.Cells.Item(195, 3).Value = '16:00'
$ws.Cells.Item(195, 4).Value = 'Bathroom'
$ws.Cells.Item(195, 5).Value = '22.8C'
$ws.Cells.Item(195, 6).Value = 'Active'

$ws.Cells.Item(196, 1).Value = '2026-01-28'
$ws.Cells.Item(196, 2).Value = '16:26:51'
$ws.Cells.Item(196, 3).Value = '16:00'
$ws.Cells.Item(196, 4).Value = 'Bathroom'
$ws.Cells.Item(196, 5).Value = '22.7C'
$ws.Cells.Item(196, 6).Value = 'Active'

$ws.Cells.Item(197, 1).Value = '2026-01-28'
$ws.Cells.Item(197, 2).Value = '16:26:55'
$ws.Cells.Item(197, 3).Value = '16:00'
$ws.Cells.Item(197, 4).Value = 'Bathroom'
$ws.Cells.Item(197, 5).Value = '22.7C'
$ws.Cells.Item(197, 6).Value = 'Active'

$ws.Range("A185:A197").Style = "Normal"
